$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.833.02'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '1.886.84'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'0.7494"
$ws.Range("E5").Value = '  -3.17%  '

$ws.Range("D6").Value = "'242.07"
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D8").Value = "'0.3124"
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = "'25.24"
$ws.Range("E9").Value = '  -2.29%  '

$ws.Range("E10").Value = '  -3.49%  '

$ws.Range("D11").Value = "'0.08505"
$ws.Range("E11").Value = '  +5.37%  '

$ws.Range("D12").Value = "'0.7588"
$ws.Range("E12").Value = '  -1.76%  '

$ws.Range("D13").Value = '1.886.37'
$ws.Range("E13").Value = '  -1.98%  '

$ws.Range("D14").Value = "'5.359"
$ws.Range("E14").Value = '  -2.68%  '

$ws.Range("D15").Value = "'93.26"
$ws.Range("E15").Value = '  -1.11%  '

$ws.Range("D16").Value = "'6.125"
$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("D17").Value = '29.851.79'
$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("E18").Value = '  -2.20%  '

$ws.Range("D19").Value = "'242.88"
$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").Value = "'0.000007833"
$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.140.25'
$ws.Range("E22").Value = '  -3.54%  '

$ws.Range("D23").Value = "'7.989"
$ws.Range("E23").Value = '  -1.95%  '

$ws.Range("D24").Value = "'0.9992"
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").Value = "'0.1582"
$ws.Range("E25").Value = '  +0.16%  '

$ws.Range("D26").Value = "'9.359"
$ws.Range("E26").Value = '  -1.16%  '

$ws.Range("D27").Value = "'163.04"
$ws.Range("E27").Value = '  -0.12%  '

$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").Value = "'1.474"
$ws.Range("E30").Value = '  +3.38%  '

$ws.Range("D31").Value = "'1.532"
$ws.Range("E31").Value = '  -0.72%  '

$ws.Range("D32").Value = "'4.504"
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("D33").Value = "'4.149"
$ws.Range("E33").Value = '  +1.98%  '

$ws.Range("D34").Value = "'0.05419"
$ws.Range("E34").Value = '  -2.82%  '

$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("D36").Value = "'0.7518"
$ws.Range("E36").Value = '  -0.22%  '

$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").Value = "'2.710"

$ws.Range("D39").Value = "'0.01943"
$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").Value = "'2.772"
$ws.Range("E40").Value = '  -0.73%  '

$ws.Range("D41").Value = "'0.4463"
$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.101.41'
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'6.099"
$ws.Range("E43").Value = '  +1.44%  '

$ws.Range("D44").Value = "'72.46"
$ws.Range("E44").Value = '  -2.98%  '

$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("D46").Value = "'1.001"

$ws.Range("D47").Value = "'7.714"
$ws.Range("E47").Value = '  +2.31%  '

$ws.Range("D48").Value = "'102.40"
$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("D49").Value = "'1.858"
$ws.Range("E49").Value = '  -1.97%  '

$ws.Range("D50").Value = "'3.030"
$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("D51").Value = '2.037.80'
$ws.Range("E51").Value = '  -1.01%  '
